# version 1.2 from testing
# Update "loginTestData" sheet (sheet 1): refresh the success-case test rows
# with a new set of username/password/expectedResult combinations, add two
# brand-new rows, and move the active selection.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Row 6 keeps its hyperlink on column B (style "Hyperlink"), just gets new text.
$ws1.Range("A6").Value = "Fawzy"
$ws1.Range("B6").Value = "Aa123456*."
$ws1.Range("C6").Value = "success"

# Row 7 used to be an (almost) empty placeholder row with only a hyperlink
# style on B7; it now becomes a normal data row, so drop that leftover style.
$ws1.Range("B7").Value = "Aa_123456789_Aa"
$ws1.Range("B7").Style = "Normal"
$ws1.Range("A7").Value = "bola"
$ws1.Range("C7").Value = "fail"

# Row 8 gets new data.
$ws1.Range("A8").Value = "3abkareem"
$ws1.Range("B8").Value = "Aa_123456789_Aa"
$ws1.Range("C8").Value = "success"

# Row 9 gets new data.
$ws1.Range("A9").Value = "zuksh"
$ws1.Range("B9").Value = "Aa_123456789_Aa"
$ws1.Range("C9").Value = "success"

# Rows 10-11 are brand new.
$ws1.Range("A10").Value = "bolt"
$ws1.Range("B10").Value = "Aa_123456789_Aa"
$ws1.Range("C10").Value = "success"

$ws1.Range("A11").Value = "messi"
$ws1.Range("B11").Value = "Aa_123456789_Aa"
$ws1.Range("C11").Value = "success"

# The B6 hyperlink still points at mailto:Abdo@1357 even though the cell's
# text was overwritten above, so Excel records the old text as the
# hyperlink's display value.
foreach ($h in $ws1.Hyperlinks) {
    if ($h.Range.Address() -eq "`$B`$6") {
        $h.TextToDisplay = "Abdo@1357"
    }
}

# Columns B and C are now the same (merged) width.
$ws1.Columns.Item(2).ColumnWidth = 14.8
$ws1.Columns.Item(3).ColumnWidth = 14.8

# "googleLogin" sheet (sheet 2): only the active selection moved.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Select() | Out-Null

# Move the active selection on sheet 1 to A8, and re-activate sheet 1 as the
# selected tab (it must remain the workbook's active sheet/tab).
$ws1.Range("A8").Select() | Out-Null
